$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 7)
$v = $cell.Value
Write-Host $v
$cell.Value = "TESTVALUE"
$v2 = $ws.Cells.Item(2,7).Value
Write-Host $v2
